$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 131129770
$ws.Range("B7").Value = 57884
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = 'Tretåig hackspett'
$ws.Range("G7").Value = 'Picoides tridactylus'
$ws.Range("H7").Value = '(Linnaeus, 1758)'
$ws.Range("P7").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q7").Value = 410609
$ws.Range("R7").Value = 7037521
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = 'Jämtland'
$ws.Range("U7").Value = 'Åre'
$ws.Range("V7").Value = 'Jämtland'
$ws.Range("W7").Value = 'Kall'
$ws.Range("Y7").Value = "'2026-02-12"
$ws.Range("Z7").Value = '14:16'
$ws.Range("AA7").Value = "'2026-02-12"
$ws.Range("AB7").Value = '14:16'
$ws.Range("AD7").Value = $False
$ws.Range("AE7").Value = $False
$ws.Range("AG7").Value = $False
$ws.Range("AW7").Value = 'Fabian Pettersson'
$ws.Range("AX7").Value = 'Fabian Pettersson'

# Row 8
$ws.Range("A8").Value = 131130868
$ws.Range("B8").Value = 57884
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = 'Tretåig hackspett'
$ws.Range("G8").Value = 'Picoides tridactylus'
$ws.Range("H8").Value = '(Linnaeus, 1758)'
$ws.Range("P8").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q8").Value = 410593
$ws.Range("R8").Value = 7037644
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 'Jämtland'
$ws.Range("U8").Value = 'Åre'
$ws.Range("V8").Value = 'Jämtland'
$ws.Range("W8").Value = 'Kall'
$ws.Range("Y8").Value = "'2026-02-06"
$ws.Range("Z8").Value = '15:36'
$ws.Range("AA8").Value = "'2026-02-06"
$ws.Range("AB8").Value = '15:36'
$ws.Range("AD8").Value = $False
$ws.Range("AE8").Value = $False
$ws.Range("AG8").Value = $False
$ws.Range("AW8").Value = 'Fabian Pettersson'
$ws.Range("AX8").Value = 'Fabian Pettersson'

# Row 9
$ws.Range("A9").Value = 131130889
$ws.Range("B9").Value = 58520
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 102125
$ws.Range("F9").Value = 'Tallbit'
$ws.Range("G9").Value = 'Pinicola enucleator'
$ws.Range("H9").Value = '(Linnaeus, 1758)'
$ws.Range("M9").Value = 'födosökande'
$ws.Range("P9").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q9").Value = 410590
$ws.Range("R9").Value = 7037649
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Jämtland'
$ws.Range("U9").Value = 'Åre'
$ws.Range("V9").Value = 'Jämtland'
$ws.Range("W9").Value = 'Kall'
$ws.Range("Y9").Value = "'2026-02-06"
$ws.Range("Z9").Value = '15:36'
$ws.Range("AA9").Value = "'2026-02-06"
$ws.Range("AB9").Value = '15:36'
$ws.Range("AD9").Value = $False
$ws.Range("AE9").Value = $False
$ws.Range("AG9").Value = $False
$ws.Range("AW9").Value = 'Fabian Pettersson'
$ws.Range("AX9").Value = 'Fabian Pettersson'

# Row 10
$ws.Range("A10").Value = 131129861
$ws.Range("B10").Value = 58520
$ws.Range("D10").Value = 'VU'
$ws.Range("E10").Value = 102125
$ws.Range("F10").Value = 'Tallbit'
$ws.Range("G10").Value = 'Pinicola enucleator'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("M10").Value = 'födosökande'
$ws.Range("P10").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q10").Value = 410588
$ws.Range("R10").Value = 7037509
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Jämtland'
$ws.Range("U10").Value = 'Åre'
$ws.Range("V10").Value = 'Jämtland'
$ws.Range("W10").Value = 'Kall'
$ws.Range("Y10").Value = "'2026-02-12"
$ws.Range("Z10").Value = '14:31'
$ws.Range("AA10").Value = "'2026-02-12"
$ws.Range("AB10").Value = '14:31'
$ws.Range("AD10").Value = $False
$ws.Range("AE10").Value = $False
$ws.Range("AG10").Value = $False
$ws.Range("AW10").Value = 'Fabian Pettersson'
$ws.Range("AX10").Value = 'Fabian Pettersson'

# Row 11
$ws.Range("A11").Value = 131131304
$ws.Range("B11").Value = 91828
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = 'Granticka'
$ws.Range("G11").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("P11").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q11").Value = 410603
$ws.Range("R11").Value = 7037541
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Jämtland'
$ws.Range("U11").Value = 'Åre'
$ws.Range("V11").Value = 'Jämtland'
$ws.Range("W11").Value = 'Kall'
$ws.Range("Y11").Value = "'2026-02-12"
$ws.Range("Z11").Value = '16:00'
$ws.Range("AA11").Value = "'2026-02-12"
$ws.Range("AB11").Value = '16:00'
$ws.Range("AD11").Value = $False
$ws.Range("AE11").Value = $False
$ws.Range("AG11").Value = $False
$ws.Range("AW11").Value = 'Fabian Pettersson'
$ws.Range("AX11").Value = 'Fabian Pettersson'

# Row 12
$ws.Range("A12").Value = 131130472
$ws.Range("B12").Value = 57884
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = 'Tretåig hackspett'
$ws.Range("G12").Value = 'Picoides tridactylus'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("P12").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q12").Value = 410737
$ws.Range("R12").Value = 7037762
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 'Jämtland'
$ws.Range("U12").Value = 'Åre'
$ws.Range("V12").Value = 'Jämtland'
$ws.Range("W12").Value = 'Kall'
$ws.Range("Y12").Value = "'2026-02-12"
$ws.Range("Z12").Value = '15:12'
$ws.Range("AA12").Value = "'2026-02-12"
$ws.Range("AB12").Value = '15:12'
$ws.Range("AD12").Value = $False
$ws.Range("AE12").Value = $False
$ws.Range("AG12").Value = $False
$ws.Range("AW12").Value = 'Fabian Pettersson'
$ws.Range("AX12").Value = 'Fabian Pettersson'

# Row 13
$ws.Range("A13").Value = 131130696
$ws.Range("B13").Value = 91808
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = 'Ullticka'
$ws.Range("G13").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H13").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("P13").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q13").Value = 410606
$ws.Range("R13").Value = 7037656
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 'Jämtland'
$ws.Range("U13").Value = 'Åre'
$ws.Range("V13").Value = 'Jämtland'
$ws.Range("W13").Value = 'Kall'
$ws.Range("Y13").Value = "'2026-02-12"
$ws.Range("Z13").Value = '15:35'
$ws.Range("AA13").Value = "'2026-02-12"
$ws.Range("AB13").Value = '15:35'
$ws.Range("AD13").Value = $False
$ws.Range("AE13").Value = $False
$ws.Range("AG13").Value = $False
$ws.Range("AW13").Value = 'Fabian Pettersson'
$ws.Range("AX13").Value = 'Fabian Pettersson'

# Row 14
$ws.Range("A14").Value = 131130718
$ws.Range("B14").Value = 57881
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 100049
$ws.Range("F14").Value = 'Spillkråka'
$ws.Range("G14").Value = 'Dryocopus martius'
$ws.Range("H14").Value = '(Linnaeus, 1758)'
$ws.Range("P14").Value = 'Stor-Grundsviken, Stor-Grundsviken, Jmt'
$ws.Range("Q14").Value = 410608
$ws.Range("R14").Value = 7037661
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 'Jämtland'
$ws.Range("U14").Value = 'Åre'
$ws.Range("V14").Value = 'Jämtland'
$ws.Range("W14").Value = 'Kall'
$ws.Range("Y14").Value = "'2026-02-12"
$ws.Range("Z14").Value = '15:36'
$ws.Range("AA14").Value = "'2026-02-12"
$ws.Range("AB14").Value = '15:36'
$ws.Range("AD14").Value = $False
$ws.Range("AE14").Value = $False
$ws.Range("AG14").Value = $False
$ws.Range("AW14").Value = 'Fabian Pettersson'
$ws.Range("AX14").Value = 'Fabian Pettersson'

